$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (value is not number-like, or is a column-E percent string) -
# Excel stores these natively as text, matching the source inline strings.
$plainUpdates = [ordered]@{
    "D2" = "27.678.22"
    "E2" = "  -1.43%  "
    "D3" = "1.880.46"
    "E3" = "  -1.29%  "
    "E4" = "  -0.09%  "
    "E5" = "  +1.06%  "
    "E6" = "  -0.16%  "
    "E7" = "  +2.08%  "
    "E8" = "  +0.31%  "
    "E9" = "  -5.67%  "
    "E10" = "  -3.07%  "
    "E11" = "  -1.43%  "
    "E12" = "  -0.05%  "
    "D13" = "1.871.99"
    "E13" = "  -3.15%  "
    "E14" = "  -0.90%  "
    "E15" = "  -2.37%  "
    "E16" = "  -0.01%  "
    "E17" = "  -2.50%  "
    "E18" = "  -2.15%  "
    "E19" = "  -0.28%  "
    "E20" = "  -3.16%  "
    "E21" = "  -0.13%  "
    "D22" = "27.655.93"
    "E22" = "  -1.47%  "
    "E23" = "  -3.31%  "
    "E24" = "  -1.20%  "
    "E25" = "  -0.34%  "
    "D26" = "2.097.51"
    "E26" = "  -2.78%  "
    "E27" = "  +0.45%  "
    "E28" = "  +1.55%  "
    "E29" = "  -1.14%  "
    "E31" = "  -1.04%  "
    "E32" = "  -0.79%  "
    "E33" = "  -1.25%  "
    "E34" = "  +0.71%  "
    "E35" = "  -0.24%  "
    "E36" = "  -3.51%  "
    "E37" = "  -0.36%  "
    "E38" = "  -1.24%  "
    "E39" = "  -2.86%  "
    "E40" = "  -5.43%  "
    "E41" = "  -2.11%  "
    "E42" = "  -0.15%  "
    "E43" = "  -0.13%  "
    "E44" = "  -4.52%  "
    "E45" = "  -2.32%  "
    "E46" = "  -4.50%  "
    "E47" = "  -4.14%  "
    "E48" = "  -0.61%  "
    "E49" = "  -3.12%  "
    "E50" = "  -1.12%  "
    "E51" = "  -0.41%  "
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Column-D price strings that look numeric (e.g. "330.68") must stay TEXT, exactly
# like the source workbook (t="inlineStr"). Assigning a bare numeric-looking string
# would auto-convert the cell to a number, so prefix with an apostrophe to force
# text entry, then restore the cell's original (unstyled) appearance so no stray
# "text" number-format/style gets attached to the cell.
$quotedUpdates = [ordered]@{
    "D5" = "330.68"
    "D7" = "0.4726"
    "D8" = "0.3986"
    "D9" = "49.19"
    "D10" = "0.08077"
    "D11" = "1.027"
    "D12" = "21.87"
    "D14" = "5.966"
    "D16" = "1.003"
    "D17" = "87.09"
    "D18" = "0.00001039"
    "D20" = "17.26"
    "D25" = "2.302"
    "D27" = "154.42"
    "D28" = "20.31"
    "D29" = "2.105"
    "D30" = "5.599"
    "D31" = "122.75"
    "D32" = "0.09513"
    "D33" = "0.9565"
    "D34" = "1.479"
    "D35" = "3.613"
    "D36" = "5.313"
    "D37" = "0.06119"
    "D39" = "1.224"
    "D40" = "8.250"
    "D41" = "0.6001"
    "D43" = "0.1900"
    "D44" = "10.36"
    "D45" = "0.5711"
    "D46" = "1.245"
    "D48" = "3.413"
    "D49" = "1.942"
    "D50" = "0.06821"
    "D51" = "110.18"
}

foreach ($ref in $quotedUpdates.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.Value = "'" + $quotedUpdates[$ref]
    $cell.Style = $origStyle
}
